$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision"" regarding the movie to be shown on Friday.`n"
$ws.Cells.Item(3, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" to be shown on Friday.`n"
$ws.Cells.Item(4, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Cells.Item(4, 4).Value = "no_decision, "
$ws.Cells.Item(5, 3).Value = "MSG: None`n`nMSG: No decision was made about which movie to show on Friday.`n"
$ws.Cells.Item(6, 3).Value = "MSG: None`n`nMSG: The decision was made that there is no consensus regarding which movie to show on Friday.`n"
$ws.Cells.Item(7, 3).Value = "MSG: None`n`nMSG: The decision has been successfully recorded, indicating that ""Barbie"" will be shown on Friday.`n"
$ws.Cells.Item(8, 3).Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for the upcoming event.`n"
$ws.Cells.Item(9, 3).Value = "MSG: None`n`nMSG: The decision was made to not acquire any movie rights, as there was no consensus on what movie to show on Friday.`n"
$ws.Cells.Item(10, 3).Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights to both movies for showing.`n"
$ws.Cells.Item(10, 4).Value = "both_movies, "
$ws.Cells.Item(11, 3).Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to show on Friday.`n"
$ws.Cells.Item(12, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights to ""Oppenheimer"" will be acquired for Friday's showing.`n"
$ws.Cells.Item(13, 3).Value = "MSG: None`n`nMSG: The function has executed, indicating that no decision was made regarding which movie to show on Friday.`n"
$ws.Cells.Item(14, 3).Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has not been made, and therefore, no movie rights will be acquired.`n"
$ws.Cells.Item(15, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" for the upcoming assembly.`n"
$ws.Cells.Item(16, 3).Value = "MSG: None`n`nMSG: The decision has been recorded. The selected movie to acquire rights for is ""Barbie.""`n"
$ws.Cells.Item(17, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Cells.Item(18, 3).Value = "MSG: None`n`nMSG: No decision was made regarding the movie to be shown on Friday.`n"
$ws.Cells.Item(19, 3).Value = "MSG: None`n`nMSG: I have recorded the decision as ""no decision"" about the movie to be shown on Friday.`n"
$ws.Cells.Item(20, 3).Value = "MSG: None`n`nMSG: None`n`nMSG: No movie was selected for the assembly on Friday, and thus I have concluded the decision process without an agreement.`n"
$ws.Cells.Item(20, 4).Value = "no_decision, , no_decision, "
$ws.Cells.Item(21, 3).Value = "MSG: None`n`nMSG: The decision resulting from the discussion indicates that no agreement was reached regarding which movie to show on Friday.`n"
$ws.Cells.Item(22, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be selected for Friday.`n"
$ws.Cells.Item(23, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Oppenheimer.""`n"
$ws.Cells.Item(24, 3).Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Cells.Item(26, 3).Value = "MSG: None`n`nMSG: The decision process has concluded without a finalized choice for Friday’s movie.`n"
$ws.Cells.Item(27, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" as the movie to be shown on Friday.`n"
$ws.Cells.Item(28, 3).Value = "MSG: None`n`nMSG: The decision concluded with no movie selected for Friday.`n"
$ws.Cells.Item(29, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie.""`n"
$ws.Cells.Item(30, 3).Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for ""Barbie"" as the movie to be shown on Friday.`n"
$ws.Cells.Item(31, 3).Value = "MSG: None`n`nMSG: No decision was made regarding the movie to show on Friday.`n"
$ws.Cells.Item(32, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached regarding the movie selection for Friday.`n"
$ws.Cells.Item(33, 3).Value = "MSG: None`n`nMSG: The decision to select a movie for Friday could not be reached, and thus no movie will be acquired.`n"
$ws.Cells.Item(34, 3).Value = "MSG: None`n`nMSG: The conversation ended without a decision about what movie to play on Friday.`n"
$ws.Cells.Item(35, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to ""Oppenheimer"" for the movie to be shown on Friday.`n"
$ws.Cells.Item(36, 3).Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for ""Barbie.""`n"
$ws.Cells.Item(37, 3).Value = "MSG: None`n`nMSG: The decision-making process concluded without a specific movie being selected for Friday.`n"
$ws.Cells.Item(38, 3).Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday was not made.`n"
$ws.Cells.Item(39, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" to be shown on Friday.`n"
$ws.Cells.Item(40, 3).Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Cells.Item(41, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached regarding the movie to be shown on Friday.`n"
$ws.Cells.Item(43, 3).Value = "MSG: None`n`nMSG: The decision to acquire the rights for ""Barbie"" has been recorded successfully.`n"
$ws.Cells.Item(44, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that there was no agreement on what movie to show on Friday.`n"
$ws.Cells.Item(45, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie.""`n"
$ws.Cells.Item(46, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Cells.Item(47, 3).Value = "MSG: None`n`nMSG: The decision results in no movie being selected for Friday.`n"
$ws.Cells.Item(48, 3).Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday resulted in no agreement.`n"
$ws.Cells.Item(49, 3).Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to ""Barbie.""`n"
$ws.Cells.Item(50, 3).Value = "MSG: None`n`nMSG: The decision has been recorded successfully, and ""Barbie"" was selected to be shown on Friday.`n"
$ws.Cells.Item(51, 3).Value = "MSG: None`n`nMSG: The decision process has concluded without an agreement on which movie to show on Friday.`n"
$ws.Cells.Item(52, 3).Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding which movie to show on Friday. Therefore, no movie rights will be acquired.`n"
$ws.Cells.Item(53, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision.""`n"
$ws.Cells.Item(53, 4).Value = "no_decision, "
$ws.Cells.Item(54, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Cells.Item(55, 3).Value = "MSG: None`n`nMSG: The decision has been recorded for acquiring the rights to ""Barbie.""`n"
$ws.Cells.Item(56, 3).Value = "MSG: None`n`nMSG: The rights for both movies have been acquired successfully.`n"
$ws.Cells.Item(57, 3).Value = "MSG: None`n`nMSG: I have recorded the decision as no decision regarding the movie to show on Friday.`n"
$ws.Cells.Item(58, 3).Value = "MSG: None`n`nMSG: The decision regarding Friday's movie ended without a clear choice, so I have recorded that as no decision made.`n"
$ws.Cells.Item(59, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be acquired at this time.`n"
$ws.Cells.Item(60, 3).Value = "MSG: None`n`nMSG: The decision reflects that no specific movie was chosen for Friday, resulting in no acquisition of rights at this time.`n"
$ws.Cells.Item(61, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday.`n"
$ws.Cells.Item(62, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision,"" indicating that no movie was selected to show on Friday.`n"
$ws.Cells.Item(63, 3).Value = "MSG: None`n`nMSG: No decision was made about the movie to be shown on Friday.`n"
$ws.Cells.Item(64, 3).Value = "MSG: None`n`nMSG: The decision has been recorded, and there is no consensus on which movie to show on Friday.`n"
$ws.Cells.Item(65, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as there being no choice of a movie.`n"
$ws.Cells.Item(66, 3).Value = "MSG: None`n`nMSG: The decision has been recorded as having no definitive choice for the movie to be shown on Friday.`n"
$ws.Cells.Item(67, 3).Value = "MSG: None`n`nMSG: The committee did not come to a decision about which movie to show on Friday.`n"
$ws.Cells.Item(68, 3).Value = "MSG: None`n`nMSG: The decision has been recorded to acquire rights for ""Oppenheimer.""`n"
$ws.Cells.Item(69, 3).Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights to both movies.`n"
$ws.Cells.Item(69, 4).Value = "both_movies, "
